# Updates cryptos list figures (Price / Volume(1h)) per the Oct 21 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '67.863.50'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -1.22%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.679.08'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -1.10%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.03%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '600.21'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +0.06%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '167.53'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +2.76%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  -0.01%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.547'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  +0.46%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '2.678.98'; ForceText = $false },
    @{ Cell = 'E9'; Value = '  -1.11%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  +3.29%  '; ForceText = $false },
    @{ Cell = 'E11'; Value = '  +1.26%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  +0.63%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '5.24'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -1.64%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '27.99'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -1.51%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '3.164.25'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -1.43%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '0.0000186'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -1.77%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '67.729.90'; ForceText = $false },
    @{ Cell = 'E17'; Value = '  -1.43%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '2.678.31'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  -1.82%  '; ForceText = $false },
    @{ Cell = 'E19'; Value = '  -0.76%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '7.81'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +1.49%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '364.93'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -0.16%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  -3.05%  '; ForceText = $false },
    @{ Cell = 'E23'; Value = '  -1.44%  '; ForceText = $false },
    @{ Cell = 'E24'; Value = '  -3.34%  '; ForceText = $false },
    @{ Cell = 'E25'; Value = '  +0.07%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '71.09'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -3.94%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '10.22'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  +3.05%  '; ForceText = $false },
    @{ Cell = 'E28'; Value = '  -0.91%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  -2.15%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  -0.26%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '560.18'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -6.10%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '8.06'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -2.90%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  -3.04%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  -0.72%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  -1.01%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '1.00'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  +0.04%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  -4.29%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '19.58'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -1.65%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '155.68'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -3.09%  '; ForceText = $false },
    @{ Cell = 'E40'; Value = '  -1.68%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '5.33'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -1.93%  '; ForceText = $false },
    @{ Cell = 'E42'; Value = '  -3.91%  '; ForceText = $false },
    @{ Cell = 'E43'; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = 'E44'; Value = '  -5.76%  '; ForceText = $false },
    @{ Cell = 'E45'; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '40.30'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -1.03%  '; ForceText = $false },
    @{ Cell = 'E47'; Value = '  -4.51%  '; ForceText = $false },
    @{ Cell = 'E48'; Value = '  -2.44%  '; ForceText = $false },
    @{ Cell = 'E49'; Value = '  -2.61%  '; ForceText = $false },
    @{ Cell = 'E50'; Value = '  -1.67%  '; ForceText = $false },
    @{ Cell = 'E51'; Value = '  -3.06%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Plain decimal-looking text (e.g. "600.21", "1.00") would otherwise be
        # auto-coerced to a number by Excel, losing the exact printed text -
        # force text format, assign, then drop back to the default style so no
        # stray formatting is left behind (matches the source cells' style 0).
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
